# Daily attendance processing - swap the order of the two recorder
# names/emails listed in the "Recorded By" column (column G) wherever the
# value consists of exactly two comma-separated entries with "System" or
# "admin@admin.com" listed before "dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Column G is "Recorded By"
$col = 7

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "admin@admin.com, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, admin@admin.com"
    }
}
